$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) stays the same content, rows 2-7 get new word/theme pairs.
$ws.Range("A1").Value = "Palavra_secreta"
$ws.Range("B1").Value = "Tema"

$ws.Range("A2").Value = "BMW"
$ws.Range("B2").Value = "Carro de Luxo"

$ws.Range("A3").Value = "Violao"
$ws.Range("B3").Value = "Instrumento"

$ws.Range("A4").Value = "Interativo"
$ws.Range("B4").Value = "Jogo da forca é ?"

$ws.Range("A5").Value = "Importante"
$ws.Range("B5").Value = "Ingles é ?"

$ws.Range("A6").Value = "Gabriel"
$ws.Range("B6").Value = "Nome do Criador desse jogo?"

$ws.Range("A7").Value = "Python"
$ws.Range("B7").Value = "Esse jogo foi feito no ?"

# Update the active selection/view to A7 (matches the saved workbook view).
$null = $ws.Range("A7").Select()
